$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add")

# Insert 3 new rows before row 29 (shifts existing rows 29+ down to 32+)
$ws.Rows("29:31").Insert()

# Populate the new rows with the readFile function documentation
$ws.Cells.Item(29, 1).Value = "roosterjs-editor-dom"
$ws.Cells.Item(29, 2).Value = "function"
$ws.Cells.Item(29, 3).Value = "readFile"
$ws.Cells.Item(29, 4).Value = "void"

$ws.Cells.Item(30, 1).Value = "roosterjs-editor-dom"
$ws.Cells.Item(30, 2).Value = "funciton.param"
$ws.Cells.Item(30, 3).Value = "readFile.file"
$ws.Cells.Item(30, 4).Value = "File"

$ws.Cells.Item(31, 1).Value = "roosterjs-editor-dom"
$ws.Cells.Item(31, 2).Value = "function.param"
$ws.Cells.Item(31, 3).Value = "readFile.callback"
$ws.Cells.Item(31, 4).Value = "function"

# Update the view: clear the frozen top-left scroll position and move the
# active selection to B13
$ws.Activate() | Out-Null
$ws.Range("B13").Select() | Out-Null
